$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# "Analysis Completion Date"  ->  three runs: "Analysis" / " End" / " Date"
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Analysis Completion Date")
if (-not $found1) {
    throw "Could not find 'Analysis Completion Date'"
}

$xml1 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:p w14:paraId="65E26314" w14:textId="77777777" w:rsidR="00033224" w:rsidRPr="00F46DE3" w:rsidRDefault="00033224" w:rsidP="00661B66"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="000000"/></w:rPr></w:pPr><w:r w:rsidRPr="00F46DE3"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="000000"/></w:rPr><w:t>Analysis</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> End</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Date</w:t></w:r></w:p></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# "Analysis Completion Time"  ->  three runs: "Analysis " / "End " / "Time"
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Analysis Completion Time")
if (-not $found2) {
    throw "Could not find 'Analysis Completion Time'"
}

$xml2 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:p w14:paraId="4ABB86FC" w14:textId="77777777" w:rsidR="00033224" w:rsidRPr="00F46DE3" w:rsidRDefault="00033224" w:rsidP="00661B66"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="000000"/></w:rPr></w:pPr><w:r w:rsidRPr="00F46DE3"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">Analysis </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">End </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="000000"/></w:rPr><w:t>Time</w:t></w:r></w:p></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng2.InsertXML($xml2)

Write-Output "Done: split 'Analysis Completion Date/Time' into Analysis/End/Date(Time) runs"
